# "run on more subjects" - an extra, never-populated column (I) that was
# sitting between the lab-timing data (A:G) and the "Participant Info"
# block (J:M) got removed, so the participant-info columns shift left by
# one (J->I, K->J, L->K, M->L) to sit right next to the timing data -
# this is what happens naturally as more subjects/columns are appended
# and the sheet gets tidied up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("In Lab")

# Delete the empty column I - this shifts J:M left by one column,
# preserving all values/styles and renumbering refs automatically.
$ws.Columns("I").Delete()

# Leave the whole (now-shifted) column I selected, matching where the
# user's cursor ended up after the delete.
$ws.Range("I1:I1048576").Select()
